$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("A1").Value = "model"
$ws.Range("B1").Value = "key function"
$ws.Range("C1").Value = "covariates"
$ws.Range("D1").Value = "Cramér von Mises p"
$ws.Range("E1").Value = "p0 ± SE"
$ws.Range("F1").Value = "AIC"
$ws.Range("G1").Value = "delta_AIC"

# Row 2 - h2 / Hazard-rate / seastate
$ws.Range("A2").Value = "h2"
$ws.Range("B2").Value = "Hazard-rate"
$ws.Range("C2").Value = "seastate"
$ws.Range("D2").Value = 0.1256
$ws.Range("E2").Value = "0.055 ± 0.0167"
$ws.Range("F2").Value = 1365.85630288711
$ws.Range("G2").Value = 0

# Row 3 - h1 / Hazard-rate / -
$ws.Range("A3").Value = "h1"
$ws.Range("B3").Value = "Hazard-rate"
$ws.Range("C3").Value = "-"
$ws.Range("D3").Value = 0.0897
$ws.Range("E3").Value = "0.1041 ± 0.0255"
$ws.Range("F3").Value = 1371.1122000113
$ws.Range("G3").Value = 5.25589712419014

# Row 4 - h3 / Hazard-rate / subj
$ws.Range("A4").Value = "h3"
$ws.Range("B4").Value = "Hazard-rate"
$ws.Range("C4").Value = "subj"
$ws.Range("D4").Value = 0.0959
$ws.Range("E4").Value = "0.1033 ± 0.0262"
$ws.Range("F4").Value = 1375.05030008077
$ws.Range("G4").Value = 9.19399719366015

# Row 5 - m1 / Half-normal / -
$ws.Range("A5").Value = "m1"
$ws.Range("B5").Value = "Half-normal"
$ws.Range("C5").Value = "-"
$ws.Range("D5").Value = 0.0055
$ws.Range("E5").Value = "0.4008 ± 0.0261"
$ws.Range("F5").Value = 1382.52054414473
$ws.Range("G5").Value = 16.66424125762

# Row 6 - m3 / Half-normal / subj
$ws.Range("A6").Value = "m3"
$ws.Range("B6").Value = "Half-normal"
$ws.Range("C6").Value = "subj"
$ws.Range("D6").Value = 0.0068
$ws.Range("E6").Value = "0.3911 ± 0.0301"
$ws.Range("F6").Value = 1383.93447284284
$ws.Range("G6").Value = 18.0781699557301

# Row 7 - m2 / Half-normal / seastate
$ws.Range("A7").Value = "m2"
$ws.Range("B7").Value = "Half-normal"
$ws.Range("C7").Value = "seastate"
$ws.Range("D7").Value = 0.0076
$ws.Range("E7").Value = "0.3933 ± 0.0279"
$ws.Range("F7").Value = 1385.14816715735
$ws.Range("G7").Value = 19.29186427024
